$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164; this pushes the existing row 164
# (and everything below it) down by one row, matching the diff which
# shows every row from 164..260 "becoming" the row below it (165..261)
# while a brand new row appears at 164.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with its data.
$ws.Range("A164").Value = 6
$ws.Range("B164").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C164").Value = "Metropolitana"
$ws.Range("D164").Value = 44438
$ws.Range("D164").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E164").Value = 13
$ws.Range("F164").Value = 100112039
$ws.Range("G164").Value = "Ciboulette"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Segunda"
$ws.Range("J164").Value = 230
$ws.Range("K164").Value = 1500
$ws.Range("L164").Value = 1500
$ws.Range("M164").Value = 1500
$ws.Range("N164").Value = "`$/docena de atados"
$ws.Range("O164").Value = "Región Metropolitana"
$ws.Range("P164").Value = 500
$ws.Range("Q164").Value = 3
$ws.Range("R164").Value = "Hortaliza"
